$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Population")

# Shift Z,AA,AB down-corrected: for r=50 down to 3, after[r] = before[r-1]
# Process from bottom (r=50) to top (r=3) so we don't clobber source values before reading them.
for ($r = 50; $r -ge 3; $r--) {
    $srcRow = $r - 1
    $z = $ws.Cells.Item($srcRow, 26).Value2
    $aa = $ws.Cells.Item($srcRow, 27).Value2
    $ab = $ws.Cells.Item($srcRow, 28).Value2
    $ws.Cells.Item($r, 26).Value2 = $z
    $ws.Cells.Item($r, 27).Value2 = $aa
    $ws.Cells.Item($r, 28).Value2 = $ab
}

# Row 2 gets the Y2 (2016) value repeated for Z/AA/AB (per source fix)
$y2 = $ws.Cells.Item(2, 25).Value2
$ws.Cells.Item(2, 26).Value2 = $y2
$ws.Cells.Item(2, 27).Value2 = $y2
$ws.Cells.Item(2, 28).Value2 = $y2
